$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# KPIs sheet: refresh top-line aggregate numbers
# ---------------------------------------------------------------------------
$wsKpis = $wb.Worksheets.Item("KPIs")
$wsKpis.Range("A2").Value = 7416.700000000001
$wsKpis.Range("E2").Value = 154.5145833333333
$wsKpis.Range("F2").Value = 0.8152173913043478

# ---------------------------------------------------------------------------
# Ventes Mensuelles sheet: updated monthly totals / quantities
# ---------------------------------------------------------------------------
$wsVentes = $wb.Worksheets.Item("Ventes Mensuelles")
$wsVentes.Range("C2").Value = 582.3
$wsVentes.Range("E2").Value = 36

$wsVentes.Range("C3").Value = 336.2
$wsVentes.Range("E3").Value = 25

$wsVentes.Range("C4").Value = 1398.7
$wsVentes.Range("E4").Value = 93

$wsVentes.Range("C5").Value = 2607.4
$wsVentes.Range("E5").Value = 180

$wsVentes.Range("C6").Value = 1087.7
$wsVentes.Range("E6").Value = 73

$wsVentes.Range("C7").Value = 1404.4
$wsVentes.Range("E7").Value = 101

# ---------------------------------------------------------------------------
# Par Catégorie sheet: updated category totals
# ---------------------------------------------------------------------------
$wsCat = $wb.Worksheets.Item("Par Catégorie")
$wsCat.Range("B2").Value = 4418.099999999999
$wsCat.Range("D2").Value = 265

$wsCat.Range("B3").Value = 2998.6
$wsCat.Range("C3").Value = 29
$wsCat.Range("D3").Value = 243

# ---------------------------------------------------------------------------
# Top Produits sheet: updated product totals
# ---------------------------------------------------------------------------
$wsTop = $wb.Worksheets.Item("Top Produits")
$wsTop.Range("B2").Value = 4418.099999999999
$wsTop.Range("C2").Value = 265

$wsTop.Range("B3").Value = 1611.5
$wsTop.Range("C3").Value = 175
$wsTop.Range("D3").Value = 29

$wsTop.Range("B4").Value = 1387.1
$wsTop.Range("C4").Value = 68
$wsTop.Range("D4").Value = 15

# ---------------------------------------------------------------------------
# Par Pays sheet: updated country total
# ---------------------------------------------------------------------------
$wsPays = $wb.Worksheets.Item("Par Pays")
$wsPays.Range("B2").Value = 7416.7

# ---------------------------------------------------------------------------
# Employés sheet: re-ranked employees with refreshed totals
# ---------------------------------------------------------------------------
$wsEmp = $wb.Worksheets.Item("Employés")

$wsEmp.Range("A2").Value = "Nancy Freehafer"
$wsEmp.Range("B2").Value = 2033.4
$wsEmp.Range("C2").Value = 12
$wsEmp.Range("D2").Value = 6

$wsEmp.Range("A3").Value = "Anne Hellung-Larsen"
$wsEmp.Range("B3").Value = 1410.1
$wsEmp.Range("C3").Value = 10
$wsEmp.Range("D3").Value = 5

$wsEmp.Range("B4").Value = 1354.6

$wsEmp.Range("B5").Value = 814.3

$wsEmp.Range("A6").Value = "Robert Zare"
$wsEmp.Range("B6").Value = 486
$wsEmp.Range("C6").Value = 2
$wsEmp.Range("D6").Value = 1

$wsEmp.Range("B7").Value = 465.8

$wsEmp.Range("A8").Value = "Michael Neipper"
$wsEmp.Range("B8").Value = 449.1
$wsEmp.Range("C8").Value = 4
$wsEmp.Range("D8").Value = 2

$wsEmp.Range("B9").Value = 403.4
